# Add "Italy" and "Spain" market test-data sheets (copies of the
# existing "Norway" template sheet), matching the structure/styling
# of the other country sheets already in the workbook.

$wb = $excel.ActiveWorkbook
$norway = $wb.Worksheets.Item("Norway")

# --- Create "Spain" first (copied right after Norway) -------------------
$norway.Copy($null, $norway)
$spainTmp = $wb.Worksheets.Item($wb.Worksheets.Count)
$spainTmp.Name = "Spain"
$spainTmp.Range("B2").Value = "Spain Market"
$spainTmp.Range("B4").Value = "NGC-3442/T2131"

# --- Create "Italy" next (copied right after Spain) ----------------------
$norway.Copy($null, $spainTmp)
$italyTmp = $wb.Worksheets.Item($wb.Worksheets.Count)
$italyTmp.Name = "Italy"
$italyTmp.Range("B2").Value = "Italy Market"
$italyTmp.Range("B4").Value = "NGC-3443/T1973"

# --- Reorder so "Italy" sits before "Spain" -------------------------------
$wb.Worksheets.Item("Italy").Move($null, $norway)

# Re-fetch worksheet references by name (positional refs go stale after Move)
$italyWs = $wb.Worksheets.Item("Italy")
$spainWs = $wb.Worksheets.Item("Spain")

# --- Match per-sheet selection state --------------------------------------
$spainWs.Activate()
$spainWs.Range("A1:XFD1048576").Select()

$italyWs.Activate()
$italyWs.Range("B4").Select()
